# Apply the edit described by the diff: the run containing
# "Connection java file's primary job is to connect the client to the
# server, ... is sent through this file's functions." is rewritten so
# it reads "...is to streamline communication between ds_client_test.java
# and the server, by having dedicated functions to handle sending and
# receiving messages. Any messages sent from the client to the server,
# are sent through this file's functions."

$d = $word.ActiveDocument

$original = "Connection java file’s primary job is to connect the client to the server, and create the ‘handshake’ between both sides, by sending messages to server, and receive and read the messages from the Server. Any messages sent from the client to the server, is sent through this file’s functions."

# Locate the exact sentence (it is unique in the document).
$range = $d.Content
$found = $range.Find.Execute($original, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target sentence"
}

$startPos = $range.Start

# Clear the old sentence.
$range.Text = ""

# Rebuild the sentence piece by piece, in the same order as the new
# runs from the diff, by repeatedly inserting after a collapsed range
# that walks forward through the freshly inserted text.
$pieces = @(
    "Connection java file’s primary job ",
    "is to streamline communication between ds_client_test.java and the server",
    ", by ",
    "having dedicated functions to handle sending and receiving messages",
    ". Any messages sent from the client to the server, ",
    "are",
    " sent through this file’s functions."
)

$cursor = $startPos
foreach ($piece in $pieces) {
    $insertionRange = $d.Range($cursor, $cursor)
    $insertionRange.InsertAfter($piece)
    $cursor = $cursor + $piece.Length
}
